$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "14-06-2021 03:45"
$ws.Range("B3").Value = "hola prueba desde poooosmmaaaan"

$ws.Range("A4").Value = "14-06-2021 03:52"
$ws.Range("B4").Value = "hola prueba desde poooosmmaaaan"

$ws.Range("A5").Value = "14-06-2021 04:01"
$ws.Range("B5").Value = "hola prueba desde poooosmmaaaan"
